$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Diseases (patient-stated)
$ws.Range("B2").Value = 4.1
$ws.Range("C2").Value = 7.3
$ws.Range("D2").Value = 4.1

# Row 3: Injuries & adverse effects
$ws.Range("B3").Value = 11.8
$ws.Range("C3").Value = 18
$ws.Range("D3").Value = 10.6

# Row 4: Other
$ws.Range("C4").Value = 6.5
$ws.Range("D4").Value = 6.9

# Row 5: Symptom - Circulatory
$ws.Range("B5").Value = 7.4
$ws.Range("C5").Value = 11.2
$ws.Range("D5").Value = 9.1

# Row 6: Symptom - Digestive
$ws.Range("B6").Value = 10.4
$ws.Range("C6").Value = 12.9
$ws.Range("D6").Value = 12.1

# Row 7: Symptom - General
$ws.Range("B7").Value = 4.3
$ws.Range("C7").Value = 4.1
$ws.Range("D7").Value = 5.4

# Row 8: Symptom - Genitourinary -> Symptom - Musculoskeletal
$ws.Range("A8").Value = "Symptom " + [char]0x2013 + " Musculoskeletal"
$ws.Range("B8").Value = 2.1
$ws.Range("C8").Value = 2.9
$ws.Range("D8").Value = 1.8

# Row 9: Symptom - Nervous
$ws.Range("B9").Value = 11.3
$ws.Range("C9").Value = 9.5
$ws.Range("D9").Value = 12.8

# Row 10: Symptom - Respiratory
$ws.Range("B10").Value = 38.6
$ws.Range("C10").Value = 21.2
$ws.Range("D10").Value = 34.7

# Row 11: Symptom - Skin/Hair/Nails
$ws.Range("B11").Value = 2
$ws.Range("C11").Value = 3.4
$ws.Range("D11").Value = 1.8

# Row 12: Uncodable/Unknown
$ws.Range("B12").Value = 2
$ws.Range("C12").Value = 3
$ws.Range("D12").Value = 0.9
